$wb = $excel.ActiveWorkbook

# --- Models sheet: remove the stray empty inline-string cell at I2 ---
$modelsWs = $wb.Worksheets.Item("Models")
$modelsWs.Range("I2").ClearContents()

# --- Payouts sheet: remove the redundant "Payment Frequency" column (F) ---
# Deleting the entire column shifts "Notes" (G) left into F and
# shrinks the used range from A1:G4 to A1:F4.
$payoutsWs = $wb.Worksheets.Item("Payouts")
$payoutsWs.Range("F1").EntireColumn.Delete()
